# Scrum board: add a new task row (priority "normal") and bump the first
# task's priority to "high", per the commit:
#   "Add a new task to scrum board. CodeSharper.Core/Common/ControlFlow/
#    CommandCallControlFlow: Add an overloaded constructor ... Try to add
#    a default value of Executor"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start the new row (row 4) from a copy of row 2's formatting, so every
# column picks up the same number formats / alignment / styles used by
# the rest of the table.
$ws.Range("A2:F2").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null          # xlPasteFormats

# --- New task row content -------------------------------------------------
$ws.Range("A4").Value2 = 3
$ws.Range("B4").Value2 = "normal"

$ws.Range("C4").Value2 = "Set a default Executor to CommandCallControlFlow"
$story = $ws.Range("C4")
$story.Characters(15, 8).Font.Bold = $true               # "Executor"
$story.Characters(27, 22).Font.Bold = $true               # "CommandCallControlFlow"

# Bump the existing first task's priority to "high" (frees up its old
# "normal" value, which the new row above now reuses).
$ws.Range("B2").Value2 = "high"

$ws.Range("D4").Value2 = "In the current version there is no option to a set default executor, so every time you have to pass an executor when a CommandCallControlFlow is initialized. I should change Executor property (set branch) visibility to public and mark it with an [Inject] attribute for supporting DI. "

$ws.Range("E4").Value2 = "to-do"
$ws.Range("F4").Value2 = 42004

# Match the wrapped-text row height used for the other long rows.
$ws.Rows.Item(4).RowHeight = 75

# Leave the selection on the newly added description cell.
$ws.Range("D4").Select() | Out-Null
